# Added New Mac-Address and Document Types
# Appends 5 new device rows (ids 3000176-3000180) to Sheet1, mirroring the
# existing "Finger Print Scanner / IRIS Scanner / Web Camera / Document
# Scanner / Printer" quintuple pattern already present in the sheet (the
# "31" series), now adding a "32" series.
#
# Cells are written column-by-column (all of column B, then all of column C,
# then all of column D, ...) rather than row-by-row so that new entries land
# in the shared-strings table in the same order the source workbook shows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 157
$ids      = @(3000176, 3000177, 3000178, 3000179, 3000180)
$names    = @("Finger Print Scanner 32", "IRIS Scanner 32", "Web Camera 32", "Document Scanner 32", "Printer 32")
$macs     = @("80-75-40-E8-CA-24", "0E-1A-14-4A-6D-3A", "65-13-7F-0F-F7-53", "73-C4-DE-8E-C9-8D", "EC-74-AB-E0-0F-38")
$serials  = @("BS563Q2230824", "BS563Q2230825", "BS563Q2230826", "BS563Q2230827", "BS563Q2230828")
$dspecs   = @(165, 327, 736, 801, 920)

# Column A - id
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Range("A" + ($firstRow + $i)).Value = $ids[$i]
}

# Column B - name
for ($i = 0; $i -lt $names.Length; $i++) {
    $ws.Range("B" + ($firstRow + $i)).Value = $names[$i]
}

# Column C - mac_address
for ($i = 0; $i -lt $macs.Length; $i++) {
    $ws.Range("C" + ($firstRow + $i)).Value = $macs[$i]
}

# Column D - serial_num
for ($i = 0; $i -lt $serials.Length; $i++) {
    $ws.Range("D" + ($firstRow + $i)).Value = $serials[$i]
}

# Column E - ip_address: intentionally left blank (matches existing rows)

# Column F - dspec_id
for ($i = 0; $i -lt $dspecs.Length; $i++) {
    $ws.Range("F" + ($firstRow + $i)).Value = $dspecs[$i]
}

# Column G - lang_code
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Range("G" + ($firstRow + $i)).Value = "eng"
}

# Column H - is_active (boolean, left-aligned like the rest of the column)
$lastRow = $firstRow + $ids.Length - 1
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Range("H" + ($firstRow + $i)).Value = $true
}
$ws.Range("H" + $firstRow + ":H" + $lastRow).HorizontalAlignment = -4131

# Column I - cr_by
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Range("I" + ($firstRow + $i)).Value = "superadmin"
}

# Column J - cr_dtimes
for ($i = 0; $i -lt $ids.Length; $i++) {
    $ws.Range("J" + ($firstRow + $i)).Value = "now()"
}

# Reposition the view the way it ended up in the saved workbook: scrolled so
# row 113 is at the top, with columns K:XFD (the empty area to the right of
# the table) selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 113
$win.ScrollColumn = 1
$ws.Range("K1:XFD1048576").Select() | Out-Null
